# Auto-generated Word COM-interop script
# Applies the OOXML diff by replacing each affected paragraph's
# Range content with exact target OOXML (including w:proofErr markers)
# via Range.InsertXML, which lets us control run-splitting precisely.

$d = $word.ActiveDocument

$p3 = @'
<w:p w14:paraId="00000003" w14:textId="77777777" w:rsidR="001A4252" w:rsidRDefault="00000000"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="-708" w:hanging="285"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Persona </w:t></w:r><w:r><w:t>(</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>CodP</w:t></w:r><w:r><w:t xml:space="preserve">, Nome, Cognome, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Email</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>, Tel</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>0</w:t></w:r><w:r><w:t xml:space="preserve">, Ruolo, Referente, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PartecipaprogFin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>
'@
$d.Paragraphs(3).Range.InsertXML($p3)

$p4 = @'
<w:p w14:paraId="00000004" w14:textId="61E39E3F" w:rsidR="001A4252" w:rsidRDefault="00000000"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="-708" w:hanging="285"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Scuola</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>CodMec</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, Nome, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Prov</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CicloIstruz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, Finanziamento, TipoFin</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>0</w:t></w:r><w:r><w:t>, Collabora)</w:t></w:r></w:p>
'@
$d.Paragraphs(4).Range.InsertXML($p4)

$p5 = @'
<w:p w14:paraId="00000005" w14:textId="007DDCC9" w:rsidR="001A4252" w:rsidRDefault="00000000"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="-708" w:hanging="285"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Classe</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>CodC</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, Ordine, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TipoScuola</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00C30B90"><w:t>DocRif</w:t></w:r><w:r w:rsidR="00C30B90"><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>PERSONA</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>
'@
$d.Paragraphs(5).Range.InsertXML($p5)

$p6 = @'
<w:p w14:paraId="00000006" w14:textId="77777777" w:rsidR="001A4252" w:rsidRDefault="00000000"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="-708" w:right="-182" w:hanging="285"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Rilevazione</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>CodR</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RespIns</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>DataRil</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>DataIns</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RespRil</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ModAcquisizione</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>InfoAmb</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>INFOAMBIENTALI</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>
'@
$d.Paragraphs(6).Range.InsertXML($p6)

$p7 = @'
<w:p w14:paraId="00000007" w14:textId="77777777" w:rsidR="001A4252" w:rsidRDefault="00000000"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="-708" w:hanging="285"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>InfoAmbientali</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>CodInfo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>LargChioma</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>LungChioma</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PesoFrescChioma</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PesoSecChioma</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>AltPianta</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>LungRadici</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PesoFrescRadici</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PesoSecRadici</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>NumFiori</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>NumFrutti</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>NumFoglieDann</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SupDann</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ph</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, Umidità, Temperatura)</w:t></w:r></w:p>
'@
$d.Paragraphs(7).Range.InsertXML($p7)

$p8 = @'
<w:p w14:paraId="00000008" w14:textId="69F7A127" w:rsidR="001A4252" w:rsidRDefault="00000000"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="-708" w:hanging="285"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Dispositivo</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>IDDisp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, Tipo)</w:t></w:r></w:p>
'@
$d.Paragraphs(8).Range.InsertXML($p8)

$p9 = @'
<w:p w14:paraId="00000009" w14:textId="30AC31E9" w:rsidR="001A4252" w:rsidRDefault="00000000"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="-708" w:hanging="285"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Replica</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>NumReplica</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, Gruppo, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>DataDimora</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, Esposizione)</w:t></w:r></w:p>
'@
$d.Paragraphs(9).Range.InsertXML($p9)

$p10 = @'
<w:p w14:paraId="0000000A" w14:textId="77777777" w:rsidR="001A4252" w:rsidRDefault="00000000"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="-708" w:hanging="285"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Specie</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>NomeScientifico</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>NomeComune</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, Esposizione, Scopo, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TotRepliche</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>
'@
$d.Paragraphs(10).Range.InsertXML($p10)

$p11 = @'
<w:p w14:paraId="0000000B" w14:textId="77777777" w:rsidR="001A4252" w:rsidRDefault="00000000"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="-708" w:hanging="285"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Orto</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Nome</w:t></w:r><w:r><w:t xml:space="preserve">, Tipo, Gps, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Superf</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ContestoAmb</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, NumSensori, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TipoSensori</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>
'@
$d.Paragraphs(11).Range.InsertXML($p11)
